# "optimise 2 summon hit skill effects"
#
# Insert a new Missile row ("arrowfast" / "弓箭快") that reuses the "arrow"
# hit effect, placing it right after the existing "arrow" row (old row 5,
# Id=2). This pushes the rest of the table down by one row and bumps every
# subsequent Id by 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Insert a new sheet row at row 6 (the table's 2nd data row, right after
# the "arrow" entry). Everything below shifts down by one row.
$ws.Rows.Item(6).Insert()

# Fill in the new "arrowfast" missile entry.
$ws.Range("A6").Value = 3
$ws.Range("B6").Value = "arrowfast"
$ws.Range("C6").Value = "弓箭快"
$ws.Range("D6").Value = "arrow"
$ws.Range("E6").Value = 20
$ws.Range("F6").Value = 20
$ws.Range("G6").Value = 1
$ws.Range("H6").Value = 1

# Renumber the Id column (col A) for every row that got pushed down, so
# the sequence stays contiguous (old Id N -> new Id N+1).
for ($r = 7; $r -le 40; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 3
}

# Grow the table definition (range/autofilter) to cover the new row.
$lo.Resize($ws.Range("A3:H40"))

# Match the final selected cell recorded in the saved workbook.
$ws.Range("C5").Select()
